$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1206.625
$ws.Range("J2").Value = 2331.3333
$ws.Range("L2").Value = 2331.3333
$ws.Range("N2").Value = -2557.3333

$ws.Range("H11").Value = 235.42857
$ws.Range("I11").Value = 235.42857
$ws.Range("K11").Value = 235.42857
$ws.Range("M11").Value = -95.42857000000001

$ws.Range("H32").Value = 964.2857
$ws.Range("I32").Value = 816.8182
$ws.Range("J32").Value = 1126.5
$ws.Range("K32").Value = 816.8182
$ws.Range("L32").Value = 1126.5
$ws.Range("M32").Value = -490.8182
$ws.Range("N32").Value = -1778.5

$ws.Range("H40").Value = 1871.7142
$ws.Range("I40").Value = 1719.6666
$ws.Range("K40").Value = 1719.6666
$ws.Range("M40").Value = -1544.6666

$ws.Range("H43").Value = 1248.2
$ws.Range("I43").Value = 750
$ws.Range("J43").Value = 1995.5
$ws.Range("K43").Value = 750
$ws.Range("L43").Value = 1995.5
$ws.Range("M43").Value = -681
$ws.Range("N43").Value = -2133.5

$ws.Range("H51").Value = 4866.3335
$ws.Range("J51").Value = 4866.3335
$ws.Range("L51").Value = 4866.3335
$ws.Range("N51").Value = -5834.3335

$ws.Range("H64").Value = 5000
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 5000
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws.Range("H137").Value = 3688.25
$ws.Range("I137").Value = 2000
$ws.Range("J137").Value = 3929.4285
$ws.Range("K137").Value = 6000
$ws.Range("L137").Value = 11788.2855
$ws.Range("M137").Value = -3450
$ws.Range("N137").Value = -16888.2855

$ws.Range("H138").Value = 3242.0588
$ws.Range("I138").Value = 1706.6
$ws.Range("J138").Value = 3881.8333
$ws.Range("K138").Value = 5119.799999999999
$ws.Range("L138").Value = 11645.4999
$ws.Range("M138").Value = 20.20000000000073
$ws.Range("N138").Value = -21925.4999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1116.8572
$ws.Range("J2").Value = 1167
$ws.Range("L2").Value = 1167
$ws.Range("N2").Value = -1393

$ws.Range("H31").Value = 3950
$ws.Range("I31").Value = 3950
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 3950
$ws.Range("L31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -3656

$ws.Range("H32").Value = 2347.1365
$ws.Range("J32").Value = 284
$ws.Range("L32").Value = 284
$ws.Range("N32").Value = -858

$ws.Range("H61").Value = 5766.1
$ws.Range("I61").Value = 5523.143
$ws.Range("J61").Value = 6333
$ws.Range("K61").Value = 5523.143
$ws.Range("L61").Value = 6333
$ws.Range("M61").Value = -5311.143
$ws.Range("N61").Value = -6757

$ws.Range("H74").Value = 3811.125
$ws.Range("I74").Value = 3811.125
$ws.Range("K74").Value = 3811.125
$ws.Range("M74").Value = -2937.125

$ws.Range("H77").Value = 3811.125
$ws.Range("I77").Value = 3811.125
$ws.Range("K77").Value = 19055.625
$ws.Range("M77").Value = -14687.625

$ws.Range("H110").Value = 2115.923
$ws.Range("I110").Value = 2051.6
$ws.Range("J110").Value = 2330.3333
$ws.Range("K110").Value = 2051.6
$ws.Range("L110").Value = 2330.3333
$ws.Range("M110").Value = -6.599999999999909
$ws.Range("N110").Value = -6420.3333

$ws.Range("H116").Value = 1116.8572
$ws.Range("J116").Value = 1167
$ws.Range("L116").Value = 1167
$ws.Range("N116").Value = -5755

$ws.Range("H132").Value = 1542.88
$ws.Range("I132").Value = 1453.65
$ws.Range("J132").Value = 1899.8
$ws.Range("K132").Value = 4360.950000000001
$ws.Range("L132").Value = 5699.4
$ws.Range("M132").Value = -1830.950000000001
$ws.Range("N132").Value = -10759.4

$ws.Range("H136").Value = 5766.1
$ws.Range("I136").Value = 5523.143
$ws.Range("J136").Value = 6333
$ws.Range("K136").Value = 16569.429
$ws.Range("L136").Value = 18999
$ws.Range("M136").Value = -14019.429
$ws.Range("N136").Value = -24099

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1116.8572
$ws.Range("J3").Value = 1167
$ws.Range("L3").Value = 1167
$ws.Range("N3").Value = -1395

$ws.Range("H38").Value = 20000
$ws.Range("J38").Value = 20000
$ws.Range("L38").Value = 20000
$ws.Range("N38").Value = -20832

$ws.Range("H94").Value = 2369.6924
$ws.Range("I94").Value = 2369.6924
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 2369.6924
$ws.Range("L94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -1918.6924

$ws.Range("H134").Value = 3863.125
$ws.Range("J134").Value = 2000
$ws.Range("L134").Value = 6000
$ws.Range("N134").Value = -11070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2117.8667
$ws.Range("I31").Value = 1543.6562
$ws.Range("J31").Value = 3531.3076
$ws.Range("K31").Value = 1543.6562
$ws.Range("L31").Value = 3531.3076
$ws.Range("M31").Value = -1248.6562
$ws.Range("N31").Value = -4121.3076

$ws.Range("H34").Value = 2117.8667
$ws.Range("I34").Value = 1543.6562
$ws.Range("J34").Value = 3531.3076
$ws.Range("K34").Value = 1543.6562
$ws.Range("L34").Value = 3531.3076
$ws.Range("M34").Value = -1341.6562
$ws.Range("N34").Value = -3935.3076

$ws.Range("H45").Value = 8000
$ws.Range("I45").Value = 6000
$ws.Range("J45").Value = 10000
$ws.Range("K45").Value = 6000
$ws.Range("L45").Value = 10000
$ws.Range("M45").Value = -5407
$ws.Range("N45").Value = -11186

$ws.Range("H107").Value = 1017.9
$ws.Range("I107").Value = 811.4286
$ws.Range("K107").Value = 811.4286
$ws.Range("M107").Value = 1108.5714

$ws.Range("H132").Value = 3397.5
$ws.Range("I132").Value = 3397.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 10192.5
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -7662.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 336
$ws.Range("J7").Value = 393
$ws.Range("L7").Value = 1179
$ws.Range("N7").Value = -1403

$ws.Range("H63").Value = 1000
$ws.Range("I63").Value = 1000
$ws.Range("K63").Value = 3000
$ws.Range("M63").Value = -2251

$ws.Range("H66").Value = 1000
$ws.Range("I66").Value = 1000
$ws.Range("K66").Value = 9000
$ws.Range("M66").Value = -5256

$ws.Range("H68").Value = 1647.2
$ws.Range("J68").Value = 2000
$ws.Range("L68").Value = 6000
$ws.Range("N68").Value = -7622

$ws.Range("H71").Value = 1647.2
$ws.Range("J71").Value = 2000
$ws.Range("L71").Value = 18000
$ws.Range("N71").Value = -26112

$ws.Range("H80").Value = 13584.5
$ws.Range("I80").Value = 12865
$ws.Range("J80").Value = 13664.444
$ws.Range("K80").Value = 38595
$ws.Range("L80").Value = 40993.33199999999
$ws.Range("M80").Value = -37659
$ws.Range("N80").Value = -42865.33199999999

$ws.Range("H83").Value = 13584.5
$ws.Range("I83").Value = 12865
$ws.Range("J83").Value = 13664.444
$ws.Range("K83").Value = 115785
$ws.Range("L83").Value = 122979.996
$ws.Range("M83").Value = -111105
$ws.Range("N83").Value = -132339.996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1319.5
$ws.Range("I97").Value = 1040
$ws.Range("K97").Value = 1040
$ws.Range("M97").Value = -544

$ws.Range("H132").Value = 2174.25
$ws.Range("I132").Value = 2174.25
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6522.75
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -3992.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 995.0833
$ws.Range("I22").Value = 1065.6666
$ws.Range("J22").Value = 783.3333
$ws.Range("K22").Value = 1065.6666
$ws.Range("L22").Value = 783.3333
$ws.Range("M22").Value = -770.6666
$ws.Range("N22").Value = -1373.3333

$ws.Range("H27").Value = 995.0833
$ws.Range("I27").Value = 1065.6666
$ws.Range("J27").Value = 783.3333
$ws.Range("K27").Value = 1065.6666
$ws.Range("L27").Value = 783.3333
$ws.Range("M27").Value = -958.6666
$ws.Range("N27").Value = -997.3333

$ws.Range("H94").Value = 21500
$ws.Range("J94").Value = 40000
$ws.Range("L94").Value = 40000
$ws.Range("N94").Value = -41352

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 47500
$ws.Range("J82").Value = 47500
$ws.Range("L82").Value = 47500
$ws.Range("N82").Value = -48266

$ws.Range("H85").Value = 47500
$ws.Range("J85").Value = 47500
$ws.Range("L85").Value = 47500
$ws.Range("N85").Value = -50152
